# Rename ObjTables document/table attributes to lowerCamelCase.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$ws1.Range("A2").Value = "!!ObjTables type='Data' id='Example0'"

$ws2.Range("A1").Value = "!!ObjTables type='Data' id='Example1'"
